$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B = "harvester" - holly added S.GISH to harvester in bioSamples.
# Update every data row (2-22) in column B from "Retrofitted_1573" to "S.GISH".
$ws.Range("B2:B22").Value = "S.GISH"

# Cosmetic follow-up matching the author's interactive edit: the harvester
# column was selected/reviewed and widened slightly to fit the new text.
$ws.Columns.Item(2).ColumnWidth = 8.83
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Range("B:B").Select()
